$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Buzzer/Speaker ()"

$ws.Range("J4").Value = 120
$ws.Range("I5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0

$ws.Range("J11").Select() | Out-Null
